# Excel COM-interop edit script
# Applies planetary_positions.xlsx diff: updates birth-chart input data (Sheet 1)
# and the planetary positions table (Sheet 2).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: birth-data inputs
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet 1")

# Date of birth (kept as plain text, like the original cell)
$ws1.Range("B2").NumberFormat = "@"
$ws1.Range("B2").Value = "10/05/1999"
$ws1.Range("B2").NumberFormat = "General"

# Place of birth
$ws1.Range("B4").Value = "Bengaluru,India"

# Latitude / Longitude
$ws1.Range("B5").Value = 12.98815675
$ws1.Range("B6").Value = 77.62260003796

# The title merge band shrinks from A1:K1 to A1:G1 because the table now
# only spans columns A-G.
$ws1.Range("A1:K1").UnMerge()
$ws1.Range("A1:G1").Merge()

# Touch G10 (with a "no line" border, i.e. no visual change) so the sheet's
# used range / dimension extends to column G, matching A1:G10.
$ws1.Range("G10").Borders.LineStyle = -4142

# ---------------------------------------------------------------------------
# Sheet 2: planetary positions table
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet 2")

$ws2.Range("B3").Value = "Aries"
$ws2.Range("C3").Value = "Mars"
$ws2.Range("D3").Value = "Bharani"
$ws2.Range("E3").Value = "Venus"
$ws2.Range("F3").Value = 24.96816340468312
$ws2.Range("G3").Value = "Direct"
$ws2.Range("H3").Value = "No"
$ws2.Range("J3").Value = 2
$ws2.Range("B4").Value = "Aquarius"
$ws2.Range("C4").Value = "Saturn"
$ws2.Range("D4").Value = "Shatabhisha"
$ws2.Range("E4").Value = "Rahu"
$ws2.Range("F4").Value = 309.075052640371
$ws2.Range("G4").Value = "Direct"
$ws2.Range("H4").Value = "No"
$ws2.Range("J4").Value = 12
$ws2.Range("B5").Value = "Aries"
$ws2.Range("C5").Value = "Mars"
$ws2.Range("D5").Value = "Ashwini"
$ws2.Range("E5").Value = "Ketu"
$ws2.Range("F5").Value = 7.96109256468479
$ws2.Range("G5").Value = "Direct"
$ws2.Range("H5").Value = "No"
$ws2.Range("J5").Value = 2
$ws2.Range("B6").Value = "Gemini"
$ws2.Range("C6").Value = "Mercury"
$ws2.Range("D6").Value = "Ardra"
$ws2.Range("E6").Value = "Rahu"
$ws2.Range("F6").Value = 67.50609350739118
$ws2.Range("G6").Value = "Direct"
$ws2.Range("H6").Value = "No"
$ws2.Range("J6").Value = 4
$ws2.Range("B7").Value = "Libra"
$ws2.Range("C7").Value = "Venus"
$ws2.Range("D7").Value = "Chitra"
$ws2.Range("E7").Value = "Mars"
$ws2.Range("F7").Value = 184.8265144931989
$ws2.Range("G7").Value = "Retro"
$ws2.Range("H7").Value = "No"
$ws2.Range("J7").Value = 8
$ws2.Range("B8").Value = "Pisces"
$ws2.Range("C8").Value = "Jupiter"
$ws2.Range("D8").Value = "Revati"
$ws2.Range("E8").Value = "Mercury"
$ws2.Range("F8").Value = 356.3771724279898
$ws2.Range("G8").Value = "Direct"
$ws2.Range("H8").Value = "No"
$ws2.Range("J8").Value = 1
$ws2.Range("B9").Value = "Aries"
$ws2.Range("C9").Value = "Mars"
$ws2.Range("D9").Value = "Bharani"
$ws2.Range("E9").Value = "Venus"
$ws2.Range("F9").Value = 14.49292394624836
$ws2.Range("G9").Value = "Direct"
$ws2.Range("H9").Value = "Combust"
$ws2.Range("J9").Value = 2
$ws2.Range("B10").Value = "Capricorn"
$ws2.Range("C10").Value = "Saturn"
$ws2.Range("D10").Value = "Shravana"
$ws2.Range("E10").Value = "Moon"
$ws2.Range("F10").Value = 292.8890339183791
$ws2.Range("G10").Value = "Direct"
$ws2.Range("H10").Value = "No"
$ws2.Range("J10").Value = 11
$ws2.Range("B11").Value = "Capricorn"
$ws2.Range("C11").Value = "Saturn"
$ws2.Range("D11").Value = "Shravana"
$ws2.Range("E11").Value = "Moon"
$ws2.Range("F11").Value = 280.5216800641753
$ws2.Range("G11").Value = "Retro"
$ws2.Range("H11").Value = "No"
$ws2.Range("J11").Value = 11
$ws2.Range("B12").Value = "Scorpio"
$ws2.Range("C12").Value = "Mars"
$ws2.Range("D12").Value = "Anuradha"
$ws2.Range("E12").Value = "Saturn"
$ws2.Range("F12").Value = 225.8449774019211
$ws2.Range("G12").Value = "Retro"
$ws2.Range("H12").Value = "No"
$ws2.Range("J12").Value = 9
$ws2.Range("B13").Value = "Cancer"
$ws2.Range("C13").Value = "Moon"
$ws2.Range("D13").Value = "Ashlesha"
$ws2.Range("E13").Value = "Mercury"
$ws2.Range("F13").Value = 113.7254071935777
$ws2.Range("G13").Value = "Retro"
$ws2.Range("H13").Value = "No"
$ws2.Range("J13").Value = 5
$ws2.Range("B14").Value = "Capricorn"
$ws2.Range("C14").Value = "Saturn"
$ws2.Range("D14").Value = "Dhanishta"
$ws2.Range("E14").Value = "Mars"
$ws2.Range("F14").Value = 293.560692356362
$ws2.Range("G14").Value = "Retro"
$ws2.Range("H14").Value = "No"
$ws2.Range("J14").Value = 11
